# Apply cryptos-list price/volume refresh per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    # Force text storage so numeric-looking strings (e.g. "230.73", "0.140")
    # keep their exact original formatting/precision instead of becoming floats.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

$ws.Range("D2").Value = "44.270.57"
$ws.Range("E2").Value = "  +4.98%  "
$ws.Range("D3").Value = "2.277.80"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextCell $ws.Range("D5") "230.73"
$ws.Range("E5").Value = "  -0.60%  "
Set-TextCell $ws.Range("D6") "0.625"
$ws.Range("E6").Value = "  -0.49%  "
Set-TextCell $ws.Range("D7") "60.74"
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +4.72%  "
Set-TextCell $ws.Range("D10") "0.0941"
$ws.Range("E10").Value = "  +4.67%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "2.619.47"
$ws.Range("E12").Value = "  +1.67%  "
Set-TextCell $ws.Range("D13") "24.65"
$ws.Range("E13").Value = "  +11.37%  "
Set-TextCell $ws.Range("D14") "15.56"
$ws.Range("E14").Value = "  -1.07%  "
Set-TextCell $ws.Range("D15") "5.81"
$ws.Range("E15").Value = "  +3.78%  "
Set-TextCell $ws.Range("D16") "0.806"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "2.280.14"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "44.150.06"
$ws.Range("E18").Value = "  +4.69%  "
$ws.Range("D19").Value = "0.0₃0941"
$ws.Range("E19").Value = "  +3.89%  "
Set-TextCell $ws.Range("D20") "73.10"
$ws.Range("E20").Value = "  +1.15%  "
Set-TextCell $ws.Range("D21") "6.23"
$ws.Range("E21").Value = "  +3.26%  "
Set-TextCell $ws.Range("D22") "254.06"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +0.02%  "
Set-TextCell $ws.Range("D24") "2.58"
$ws.Range("E24").Value = "  +7.67%  "
Set-TextCell $ws.Range("D25") "2.49"
$ws.Range("E25").Value = "  +4.42%  "
Set-TextCell $ws.Range("D26") "9.83"
$ws.Range("E26").Value = "  +1.64%  "
Set-TextCell $ws.Range("D27") "171.43"
$ws.Range("E27").Value = "  +1.50%  "
Set-TextCell $ws.Range("D28") "0.140"
$ws.Range("E28").Value = "  -2.90%  "
Set-TextCell $ws.Range("D29") "20.56"
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -0.37%  "
Set-TextCell $ws.Range("D32") "0.122"
$ws.Range("E32").Value = "  -0.25%  "
Set-TextCell $ws.Range("D33") "5.04"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +1.09%  "
Set-TextCell $ws.Range("D35") "0.0654"
$ws.Range("E35").Value = "  +2.44%  "
Set-TextCell $ws.Range("D36") "6.48"
$ws.Range("E36").Value = "  -2.94%  "
Set-TextCell $ws.Range("D37") "2.38"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("E39").Value = "  +4.09%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +1.74%  "
Set-TextCell $ws.Range("D42") "0.000223"
$ws.Range("E42").Value = "  -13.29%  "
$ws.Range("D47").Value = "1.476.40"
$ws.Range("E47").Value = "  -0.46%  "
Set-TextCell $ws.Range("D48") "16.61"
$ws.Range("E48").Value = "  +0.37%  "
Set-TextCell $ws.Range("D49") "1.08"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("E51").Value = "  +5.71%  "

# Rows 43-46 were re-ranked (Cronos/FTXToken and Aave/TrustWalletToken swapped order)
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws.Range("D43") "0.0966"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell $ws.Range("D44") "4.45"
$ws.Range("E44").Value = "  -7.89%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws.Range("D45") "98.16"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Range("D46") "1.21"
$ws.Range("E46").Value = "  -1.46%  "
